$d = $word.ActiveDocument

# Locate the paragraph ending with "...presence of microcephaly at birth?"
# and insert a new bulleted list item right after it, matching the
# formatting (NoSpacing style, same numbered/bulleted list) of its
# neighbours.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*presence of microcephaly at birth?*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph"
}

$target.Range.InsertParagraphAfter()
$newPara = $target.Next()
$newPara.Range.Text = "Help with passive imputation of variables."

Write-Output "Inserted new list item."
